$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_SD"

$ws.Range("A2").Value = 71.790540540540533
$ws.Range("A3").Value = 71.621621621621628
$ws.Range("A4").Value = 71.28378378378379
$ws.Range("A5").Value = 68.75
$ws.Range("A6").Value = 69.594594594594597
$ws.Range("A7").Value = 69.932432432432435
$ws.Range("A8").Value = 71.28378378378379
$ws.Range("A9").Value = 71.452702702702695
$ws.Range("A10").Value = 71.959459459459467
$ws.Range("A11").Value = 72.128378378378372
$ws.Range("A12").Value = 68.75
$ws.Range("A13").Value = 68.75
$ws.Range("A14").Value = 69.594594594594597
$ws.Range("A15").Value = 69.425675675675677
$ws.Range("A16").Value = 69.425675675675677
$ws.Range("A17").Value = 62.162162162162161
$ws.Range("A18").Value = 67.060810810810807
$ws.Range("A19").Value = 68.581081081081081
$ws.Range("A20").Value = 71.452702702702695
$ws.Range("A21").Value = 71.959459459459467
$ws.Range("A22").Value = 71.621621621621628
$ws.Range("A23").Value = 66.722972972972968
$ws.Range("A24").Value = 66.554054054054063
$ws.Range("A25").Value = 66.554054054054063
$ws.Range("A26").Value = 71.452702702702695
$ws.Range("A27").Value = 70.101351351351354
$ws.Range("A28").Value = 70.439189189189193
$ws.Range("A29").Value = 68.412162162162161
$ws.Range("A30").Value = 67.905405405405403
$ws.Range("A31").Value = 67.736486486486484
$ws.Range("A32").Value = 72.128378378378372
$ws.Range("A33").Value = 70.439189189189193
$ws.Range("A34").Value = 71.959459459459467
$ws.Range("A35").Value = 68.074324324324323
$ws.Range("A36").Value = 68.412162162162161
$ws.Range("A37").Value = 59.29054054054054
$ws.Range("A38").Value = 68.918918918918919
$ws.Range("A39").Value = 67.060810810810807
$ws.Range("A40").Value = 68.074324324324323
$ws.Range("A41").Value = 68.581081081081081
$ws.Range("A42").Value = 68.412162162162161
$ws.Range("A43").Value = 68.75
$ws.Range("A44").Value = 68.75
$ws.Range("A45").Value = 69.425675675675677
$ws.Range("A46").Value = 70.270270270270274
$ws.Range("A47").Value = 66.554054054054063
$ws.Range("A48").Value = 65.03378378378379
$ws.Range("A49").Value = 70.101351351351354
